$wb = $excel.ActiveWorkbook
try {
  $r = $excel.Run("review.inquire.clean")
  Write-Host "OK: $r"
} catch {
  Write-Host "FAIL: $_"
}
